$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet
$lastRow = $ws.UsedRange.Rows.Count

# Swap the contents of column A and column B for every row (header + data)
for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $cellA.Value2 = $valB
    $cellB.Value2 = $valA
}
